$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2025-11-21 Friday" "2025-11-22 Saturday"

Replace-Text "54×30=" "55×16="
Replace-Text "55×91=" "16×14="
Replace-Text "68×98=" "19×12="
Replace-Text "52×81=" "55×77="
Replace-Text "43×56=" "31×42="
Replace-Text "70×91=" "73×25="
Replace-Text "41×61=" "42×49="
Replace-Text "91×75=" "37×67="
Replace-Text "16×50=" "86×60="
Replace-Text "48×93=" "25×66="
Replace-Text "93×82=" "80×19="
Replace-Text "89×54=" "19×77="
Replace-Text "66×74=" "21×65="
Replace-Text "45×53=" "30×60="
Replace-Text "17×56=" "71×95="
Replace-Text "64×64=" "14×65="
Replace-Text "77×71=" "66×63="
Replace-Text "87×77=" "85×79="
Replace-Text "42×15=" "99×55="
Replace-Text "82×45=" "49×63="
Replace-Text "12×43=" "77×72="
Replace-Text "18×47=" "47×70="
Replace-Text "42×76=" "89×18="
Replace-Text "44×80=" "27×16="
Replace-Text "40×16=" "43×30="
